# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.318.76"
$ws.Range("E2").Value = "'  -1.14%  "
$ws.Range("D3").Value = "'3.540.88"
$ws.Range("E3").Value = "'  +0.60%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("D5").Value = "'608.32"
$ws.Range("E5").Value = "'  +0.58%  "
$ws.Range("D6").Value = "'144.85"
$ws.Range("E6").Value = "'  -2.26%  "
$ws.Range("D7").Value = "'3.541.07"
$ws.Range("E7").Value = "'  +0.63%  "
$ws.Range("E8").Value = "'  -0.03%  "
$ws.Range("D9").Value = "'0.479"
$ws.Range("E9").Value = "'  -0.43%  "
$ws.Range("E10").Value = "'  -4.27%  "
$ws.Range("D11").Value = "'8.06"
$ws.Range("E11").Value = "'  +3.07%  "
$ws.Range("D12").Value = "'0.415"
$ws.Range("E12").Value = "'  -2.13%  "
$ws.Range("D13").Value = "'4.144.46"
$ws.Range("E13").Value = "'  +0.69%  "
$ws.Range("E14").Value = "'  -3.06%  "
$ws.Range("D15").Value = "'30.53"
$ws.Range("E15").Value = "'  -3.56%  "
$ws.Range("D16").Value = "'3.542.60"
$ws.Range("E16").Value = "'  +0.70%  "
$ws.Range("D17").Value = "'66.384.11"
$ws.Range("E17").Value = "'  -1.03%  "
$ws.Range("E18").Value = "'  -0.04%  "
$ws.Range("D19").Value = "'10.82"
$ws.Range("E19").Value = "'  +0.80%  "
$ws.Range("D20").Value = "'6.24"
$ws.Range("E20").Value = "'  -2.61%  "
$ws.Range("D21").Value = "'15.01"
$ws.Range("E21").Value = "'  -2.49%  "
$ws.Range("D22").Value = "'427.21"
$ws.Range("E22").Value = "'  -2.03%  "
$ws.Range("D23").Value = "'0.603"
$ws.Range("E23").Value = "'  -1.35%  "
$ws.Range("D24").Value = "'78.51"
$ws.Range("E24").Value = "'  -1.49%  "
$ws.Range("D25").Value = "'3.685.00"
$ws.Range("E25").Value = "'  +0.75%  "
$ws.Range("E26").Value = "'  -0.09%  "
$ws.Range("D27").Value = "'0.0000121"
$ws.Range("E27").Value = "'  -0.21%  "
$ws.Range("D28").Value = "'9.33"
$ws.Range("E28").Value = "'  -5.48%  "
$ws.Range("D29").Value = "'8.05"
$ws.Range("E29").Value = "'  -3.96%  "
$ws.Range("E31").Value = "'  +0.04%  "
$ws.Range("E32").Value = "'  -2.24%  "
$ws.Range("D33").Value = "'1.49"
$ws.Range("E33").Value = "'  -6.39%  "
$ws.Range("D34").Value = "'25.42"
$ws.Range("E34").Value = "'  +0.05%  "
$ws.Range("D35").Value = "'3.529.15"
$ws.Range("E35").Value = "'  +0.48%  "
$ws.Range("E36").Value = "'  -0.05%  "
$ws.Range("D37").Value = "'1.75"
$ws.Range("E37").Value = "'  -3.11%  "
$ws.Range("B38").Value = "'Aptos"
$ws.Range("C38").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "'7.85"
$ws.Range("E38").Value = "'  -2.27%  "
$ws.Range("B39").Value = "'NEARProtocol"
$ws.Range("C39").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "'5.63"
$ws.Range("E39").Value = "'  -4.54%  "
$ws.Range("E40").Value = "'  +0.00%  "
$ws.Range("D41").Value = "'170.60"
$ws.Range("E41").Value = "'  +0.76%  "
$ws.Range("D42").Value = "'0.0863"
$ws.Range("E42").Value = "'  -3.33%  "
$ws.Range("D43").Value = "'5.20"
$ws.Range("E43").Value = "'  -4.38%  "
$ws.Range("D44").Value = "'0.893"
$ws.Range("E44").Value = "'  -0.45%  "
$ws.Range("E45").Value = "'  -9.67%  "
$ws.Range("E46").Value = "'  -0.67%  "
$ws.Range("E47").Value = "'  -7.99%  "
$ws.Range("E48").Value = "'  -9.77%  "
$ws.Range("D49").Value = "'2.43"
$ws.Range("E49").Value = "'  -1.00%  "
$ws.Range("E50").Value = "'  -3.82%  "
$ws.Range("D51").Value = "'0.957"
$ws.Range("E51").Value = "'  -3.38%  "
